$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.447.71"
$ws.Range("E2").Value = "  +1.96%  "
$ws.Range("D3").Value = "1.853.75"
$ws.Range("E3").Value = "  +1.15%  "
$ws.Range("D4").Value = "'0.9998"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'245.41"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("D6").Value = "'0.6932"
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").Value = "'0.07668"
$ws.Range("E8").Value = "  -0.50%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "'0.3063"
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("D10").Value = "'23.52"
$ws.Range("E10").Value = "  +0.55%  "
$ws.Range("D11").Value = "'0.07754"
$ws.Range("D12").Value = "1.868.18"
$ws.Range("E12").Value = "  +2.02%  "
$ws.Range("D13").Value = "'5.146"
$ws.Range("D14").Value = "'0.6939"
$ws.Range("E14").Value = "  +1.81%  "
$ws.Range("D15").Value = "'91.03"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").Value = "'6.300"
$ws.Range("E16").Value = "  -2.29%  "
$ws.Range("D17").Value = "29.442.54"
$ws.Range("E17").Value = "  +2.06%  "
$ws.Range("D18").Value = "'0.000008281"
$ws.Range("E18").Value = "  -0.47%  "
$ws.Range("D19").Value = "2.098.07"
$ws.Range("E19").Value = "  +1.46%  "
$ws.Range("D20").Value = "'236.44"
$ws.Range("E20").Value = "  -2.40%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").Value = "'7.638"
$ws.Range("E23").Value = "  +2.43%  "
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("D25").Value = "'0.1480"
$ws.Range("E25").Value = "  -0.31%  "
$ws.Range("D26").Value = "'8.937"
$ws.Range("E26").Value = "  +1.57%  "
$ws.Range("D27").Value = "'159.92"
$ws.Range("E27").Value = "  +0.87%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").Value = "'1.530"
$ws.Range("E29").Value = "  -1.11%  "
$ws.Range("D30").Value = "'4.243"
$ws.Range("E30").Value = "  +0.47%  "
$ws.Range("D31").Value = "'4.135"
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("D32").Value = "'1.201"
$ws.Range("E32").Value = "  +1.08%  "
$ws.Range("D33").Value = "'0.05207"
$ws.Range("E33").Value = "  +1.92%  "
$ws.Range("D34").Value = "'0.7760"
$ws.Range("E34").Value = "  -0.19%  "
$ws.Range("D35").Value = "'1.870"
$ws.Range("E35").Value = "  +0.99%  "
$ws.Range("E36").Value = "  +0.44%  "
$ws.Range("D37").Value = "'2.681"
$ws.Range("E37").Value = "  -0.29%  "
$ws.Range("D38").Value = "1.329.51"
$ws.Range("E38").Value = "  +8.60%  "
$ws.Range("D39").Value = "'0.01866"
$ws.Range("E39").Value = "  +0.84%  "
$ws.Range("D40").Value = "'2.725"
$ws.Range("E40").Value = "  +1.20%  "
$ws.Range("D41").Value = "'0.9394"
$ws.Range("E41").Value = "  -2.19%  "
$ws.Range("D42").Value = "'106.08"
$ws.Range("E42").Value = "  -2.27%  "
$ws.Range("D43").Value = "'5.815"
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").Value = "'9.711"
$ws.Range("E45").Value = "  +0.76%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "'0.00000000124"
$ws.Range("E46").Value = "  +1.27%  "
$ws.Range("B47").Value = "RocketPoolETH"
$ws.Range("C47").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value = "1.996.68"
$ws.Range("E47").Value = "  +1.13%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.5226"
$ws.Range("E48").Value = "  +1.38%  "
$ws.Range("D49").Value = "'1.783"
$ws.Range("D50").Value = "'63.05"
$ws.Range("E50").Value = "  -2.22%  "
$ws.Range("D51").Value = "'0.05955"
$ws.Range("E51").Value = "  +0.77%  "
